# Auto-generated Excel COM-interop script applying the scheduled-runner price refresh
# to the Ravana_Profits workbook. Mirrors the upstream diff: per-row recomputation of
# market-price-derived columns (H-N) on the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").Value = ""

$ws.Range("H80").Value = 615.0909
$ws.Range("I80").Value = 672.125
$ws.Range("K80").Value = 2016.375
$ws.Range("M80").Value = -1018.375

$ws.Range("H83").Value = 615.0909
$ws.Range("I83").Value = 672.125
$ws.Range("K83").Value = 6049.125
$ws.Range("M83").Value = -1057.125

$ws.Range("H132").Value = 1179.5
$ws.Range("I132").Value = 1060.3334
$ws.Range("K132").Value = 3181.0002
$ws.Range("M132").Value = -651.0001999999999

$ws.Range("H137").Value = 4424.0835
$ws.Range("I137").Value = 3432.6667
$ws.Range("J137").Value = 5415.5
$ws.Range("K137").Value = 10298.0001
$ws.Range("L137").Value = 16246.5
$ws.Range("M137").Value = -7748.000100000001
$ws.Range("N137").Value = -21346.5

$ws.Range("H138").Value = 4431.5117
$ws.Range("I138").Value = 2991.6667
$ws.Range("K138").Value = 8975.000100000001
$ws.Range("M138").Value = -3835.000100000001

$ws.Range("H141").Value = 3349.524
$ws.Range("I141").Value = 2952.5
$ws.Range("K141").Value = 8857.5
$ws.Range("M141").Value = -3677.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 591.05884
$ws.Range("I2").Value = 608.1818
$ws.Range("J2").Value = 559.6667
$ws.Range("K2").Value = 608.1818
$ws.Range("L2").Value = 559.6667
$ws.Range("M2").Value = -495.1818
$ws.Range("N2").Value = -785.6667

$ws.Range("H32").Value = 3705.7036
$ws.Range("I32").Value = 3527.608
$ws.Range("K32").Value = 3527.608
$ws.Range("M32").Value = -3240.608

$ws.Range("H61").Value = 993.5833
$ws.Range("I61").Value = 811.1818
$ws.Range("K61").Value = 811.1818
$ws.Range("M61").Value = -599.1818

$ws.Range("H74").Value = 2455.348
$ws.Range("I74").Value = 2155.5386
$ws.Range("K74").Value = 2155.5386
$ws.Range("M74").Value = -1281.5386

$ws.Range("H77").Value = 2455.348
$ws.Range("I77").Value = 2155.5386
$ws.Range("K77").Value = 10777.693
$ws.Range("M77").Value = -6409.692999999999

$ws.Range("H116").Value = 591.05884
$ws.Range("I116").Value = 608.1818
$ws.Range("J116").Value = 559.6667
$ws.Range("K116").Value = 608.1818
$ws.Range("L116").Value = 559.6667
$ws.Range("M116").Value = 1685.8182
$ws.Range("N116").Value = -5147.6667

$ws.Range("H122").Value = 3213.1428
$ws.Range("I122").Value = 3213.1428
$ws.Range("K122").Value = 9639.428400000001
$ws.Range("M122").Value = -7189.428400000001

$ws.Range("H132").Value = 2386.4856
$ws.Range("I132").Value = 1197
$ws.Range("K132").Value = 3591
$ws.Range("M132").Value = -1061

$ws.Range("H136").Value = 993.5833
$ws.Range("I136").Value = 811.1818
$ws.Range("K136").Value = 2433.5454
$ws.Range("M136").Value = 116.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 591.05884
$ws.Range("I3").Value = 608.1818
$ws.Range("J3").Value = 559.6667
$ws.Range("K3").Value = 608.1818
$ws.Range("L3").Value = 559.6667
$ws.Range("M3").Value = -494.1818
$ws.Range("N3").Value = -787.6667

$ws.Range("H86").Value = 2959.3
$ws.Range("I86").Value = 2630.375
$ws.Range("J86").Value = 4275
$ws.Range("K86").Value = 2630.375
$ws.Range("L86").Value = 4275
$ws.Range("M86").Value = -1507.375
$ws.Range("N86").Value = -6521

$ws.Range("H89").Value = 2959.3
$ws.Range("I89").Value = 2630.375
$ws.Range("J89").Value = 4275
$ws.Range("K89").Value = 13151.875
$ws.Range("L89").Value = 21375
$ws.Range("M89").Value = -7535.875
$ws.Range("N89").Value = -32607

$ws.Range("H105").Value = 4799.875
$ws.Range("I105").Value = 4799.875
$ws.Range("K105").Value = 4799.875
$ws.Range("M105").Value = -3052.875

$ws.Range("H134").Value = 2826.9656
$ws.Range("I134").Value = 2872.3928
$ws.Range("K134").Value = 8617.178400000001
$ws.Range("M134").Value = -6082.178400000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3009
$ws.Range("J58").Value = 3009
$ws.Range("L58").Value = 3009
$ws.Range("N58").Value = -3415

$ws.Range("H132").Value = 3746.875
$ws.Range("I132").Value = 2996.3333
$ws.Range("J132").Value = 4197.2
$ws.Range("K132").Value = 8988.999899999999
$ws.Range("L132").Value = 12591.6
$ws.Range("M132").Value = -6458.999899999999
$ws.Range("N132").Value = -17651.6

$ws.Range("H134").Value = 2832.889
$ws.Range("I134").Value = 3124.5
$ws.Range("K134").Value = 9373.5
$ws.Range("M134").Value = -6838.5

$ws.Range("H136").Value = 3009
$ws.Range("J136").Value = 3009
$ws.Range("L136").Value = 9027
$ws.Range("N136").Value = -14127

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 169
$ws.Range("I17").Value = 199
$ws.Range("J17").Value = 79
$ws.Range("K17").Value = 597
$ws.Range("L17").Value = 237
$ws.Range("M17").Value = -428
$ws.Range("N17").Value = -575

$ws.Range("H46").Value = 2150.5
$ws.Range("I46").Value = 2150.5
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 6451.5
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = -6360.5
$ws.Range("M46").Value = ""

$ws.Range("H131").Value = 1329.8518
$ws.Range("I131").Value = 955.625
$ws.Range("J131").Value = 1487.421
$ws.Range("K131").Value = 2866.875
$ws.Range("L131").Value = 4462.263
$ws.Range("M131").Value = 2173.125
$ws.Range("N131").Value = -14542.263

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 97.818184
$ws.Range("I2").Value = 86.333336
$ws.Range("J2").Value = 149.5
$ws.Range("K2").Value = 86.333336
$ws.Range("L2").Value = 149.5
$ws.Range("M2").Value = 26.666664
$ws.Range("N2").Value = -375.5

$ws.Range("H80").Value = 7346.25
$ws.Range("I80").Value = 5753.4
$ws.Range("J80").Value = 10001
$ws.Range("K80").Value = 5753.4
$ws.Range("L80").Value = 10001
$ws.Range("M80").Value = -4755.4
$ws.Range("N80").Value = -11997

$ws.Range("H83").Value = 7346.25
$ws.Range("I83").Value = 5753.4
$ws.Range("J83").Value = 10001
$ws.Range("K83").Value = 28767
$ws.Range("L83").Value = 50005
$ws.Range("M83").Value = -23775
$ws.Range("N83").Value = -59989

$ws.Range("H102").Value = 4802.4
$ws.Range("I102").Value = 999.5
$ws.Range("K102").Value = 999.5
$ws.Range("M102").Value = 622.5

$ws.Range("H132").Value = 1209.9667
$ws.Range("I132").Value = 703.9231
$ws.Range("K132").Value = 2111.7693
$ws.Range("M132").Value = 418.2307000000001

$ws.Range("H139").Value = 68396
$ws.Range("J139").Value = 68396
$ws.Range("L139").Value = 68396
$ws.Range("N139").Value = -78676

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3137.5715
$ws.Range("I7").Value = 3137.5715
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 3137.5715
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = -3025.5715
$ws.Range("M7").Value = ""

$ws.Range("H43").Value = 189998.67
$ws.Range("J43").Value = 189998.67
$ws.Range("L43").Value = 189998.67
$ws.Range("N43").Value = -190384.67

$ws.Range("H46").Value = 176.92308

$ws.Range("H126").Value = 3137.5715
$ws.Range("I126").Value = 3137.5715
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 9412.7145
$ws.Range("L126").Value = 0
$ws.Range("N126").Value = -6942.7145
$ws.Range("M126").Value = ""

$ws.Range("H132").Value = 2665.5
$ws.Range("I132").Value = 1799.05
$ws.Range("K132").Value = 5397.15
$ws.Range("M132").Value = -2867.15

$ws.Range("H136").Value = 4903.857
$ws.Range("I136").Value = 4863.9
$ws.Range("J136").Value = 5003.75
$ws.Range("K136").Value = 14591.7
$ws.Range("L136").Value = 15011.25
$ws.Range("M136").Value = -12041.7
$ws.Range("N136").Value = -20111.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 11284.429
$ws.Range("I81").Value = 7397.4
$ws.Range("J81").Value = 21002
$ws.Range("K81").Value = 14794.8
$ws.Range("L81").Value = 42004
$ws.Range("M81").Value = -13733.8
$ws.Range("N81").Value = -44126

$ws.Range("H84").Value = 11284.429
$ws.Range("I84").Value = 7397.4
$ws.Range("J84").Value = 21002
$ws.Range("K84").Value = 73974
$ws.Range("L84").Value = 210020
$ws.Range("M84").Value = -68670
$ws.Range("N84").Value = -220628

$ws.Range("H122").Value = 1676.8422
$ws.Range("I122").Value = 1702.1111
$ws.Range("J122").Value = 1222
$ws.Range("K122").Value = 5106.3333
$ws.Range("L122").Value = 3666
$ws.Range("M122").Value = -2656.3333
$ws.Range("N122").Value = -8566

$ws.Range("H132").Value = 1859.6666
$ws.Range("I132").Value = 1285.125
$ws.Range("J132").Value = 3211.5293
$ws.Range("K132").Value = 3855.375
$ws.Range("L132").Value = 9634.5879
$ws.Range("M132").Value = -1325.375
$ws.Range("N132").Value = -14694.5879

$ws.Range("H136").Value = 1574.8064
$ws.Range("I136").Value = 1221.15
$ws.Range("J136").Value = 2217.818
$ws.Range("K136").Value = 3663.45
$ws.Range("L136").Value = 6653.454000000001
$ws.Range("M136").Value = -1113.45
$ws.Range("N136").Value = -11753.454
